$wb = $excel.ActiveWorkbook

# --- 1) Bump template version 1.0.3 -> 1.0.4 on the isa_template sheet ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.4"

# --- 2) Remove the "Factor [Experimental Factor Value]" / NCIT:C164386 term
#        source/accession columns (W:Y) from the annotationTable on Sample ---
$wsSample = $wb.Worksheets.Item("Sample")
$tbl = $wsSample.ListObjects.Item(1)

# Physically delete the 3 worksheet columns; this shifts everything to the
# right of them (growth protocol onward) left by 3 columns.
$wsSample.Range("W1:Y1").EntireColumn.Delete()

# Shrink the table definition to match the new, narrower extent.
$tbl.Resize($wsSample.Range("A1:AF2"))

# Re-assert the header row text so the table's column-name metadata
# re-syncs with the (now shifted) worksheet header cells.
$headers = @( `
    "Input [Source Name]", `
    "Characteristic [organism]", `
    "Term Source REF (OBI:0100026)", `
    "Term Accession Number (OBI:0100026)", `
    "Characteristic [Variety]", `
    "Term Source REF (NCIT:C62709)", `
    "Term Accession Number (NCIT:C62709)", `
    "Characteristic [age]", `
    "Term Source REF (EFO:0000246)", `
    "Term Accession Number (EFO:0000246)", `
    "Characteristic [plant structure development stage]", `
    "Term Source REF (PO:0009012)", `
    "Term Accession Number (PO:0009012)", `
    "Characteristic [Genotype]", `
    "Term Source REF (NCIT:C16631)", `
    "Term Accession Number (NCIT:C16631)", `
    "Characteristic [plant anatomical entity]", `
    "Term Source REF (PO:0025131)", `
    "Term Accession Number (PO:0025131)", `
    "Characteristic [plant material]", `
    "Term Source REF (FOODON:00004331)", `
    "Term Accession Number (FOODON:00004331)", `
    "Parameter [growth protocol]", `
    "Term Source REF (EFO:0003789)", `
    "Term Accession Number (EFO:0003789)", `
    "Parameter [sample collection protocol]", `
    "Term Source REF (EFO:0005518)", `
    "Term Accession Number (EFO:0005518)", `
    "Parameter [nucleic acid extraction protocol]", `
    "Term Source REF (EFO:0002944)", `
    "Term Accession Number (EFO:0002944)", `
    "Output [Sample Name]" `
)

for ($i = 1; $i -le $headers.Count; $i++) {
    $wsSample.Cells.Item(1, $i).Value = $headers[$i - 1]
}
